$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the value of B3 to the new shared string.
$ws.Range("B3").Value = "Спортивные снаряды"

# Move the selection / view back to B3, dropping the scrolled topLeftCell.
$ws.Range("B3").Select()
